# Append the 2025-06-20 price point to the "Gaz" and "CO2" sheets,
# matching the new row 6 added to docs/epexspot_prices.xlsx.

$wb = $excel.ActiveWorkbook

# --- "Gaz" sheet: add row 6 (Date = 2025-06-20, Last Price = 39.7) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force column A to be read as plain text so the date-like string isn't
# auto-converted into a serial date value, then reset the cell style back
# to match the other data rows (no explicit style) once the value is set.
$wsGaz.Range("A6").NumberFormat = "@"
$wsGaz.Range("A6").Value = "2025-06-20"
$wsGaz.Range("A6").Style = $wsGaz.Range("A5").Style
$wsGaz.Range("B6").Value = 39.7

# --- "CO2" sheet: add row 6 (Date = 2025-06-20, Last Price = 72.2) ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A6").NumberFormat = "@"
$wsCo2.Range("A6").Value = "2025-06-20"
$wsCo2.Range("A6").Style = $wsCo2.Range("A5").Style
$wsCo2.Range("B6").Value = 72.2
